# add properties for ELine in QStudioSCADA and QSCADARunTime.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("画面编辑器")

# Row 6 (ELine / 直线): mark as completed (已完成) -> green fill (matches style used on rows 2-5)
$ws.Range("B6").Value = "已完成"
$ws.Range("B6").Interior.Color = 5287936

# Row 7 (箭头): move to in-progress (进行中) -> yellow fill (matches style used on row 6 previously)
$ws.Range("B7").Value = "进行中"
$ws.Range("B7").Interior.Color = 65535

# Update active selection to B6
$ws.Range("B6").Select()
